# Add per-parameter-block min/max/diff summary columns (M,N,O) to Sheet1.
# Each of the 5 data blocks (rows 3-10, 18-25, 30-37, 42-49, 54-61) gets:
#   - a "min"/"max"/"diff" header on its header row (the row right above the
#     first data row), styled like the existing G:K header (style of G2/G17/...)
#   - MIN/MAX/ABS(diff) formulas on the block's first data row
#   - (for blocks 3-5 only) the M/N number-format carried down the remaining
#     data rows as empty, formatted cells (mirrors a drag-fill of the format
#     without values)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blocks = @(
    @{hdr=2;  first=3;  last=10; fillDown=$false},
    @{hdr=17; first=18; last=25; fillDown=$false},
    @{hdr=29; first=30; last=37; fillDown=$true},
    @{hdr=41; first=42; last=49; fillDown=$true},
    @{hdr=53; first=54; last=61; fillDown=$true}
)

foreach ($b in $blocks) {
    $hdr = $b.hdr
    $first = $b.first
    $last = $b.last

    # Header row: min / max / diff labels, styled like the block's G-K header cells.
    $ws.Range("M$hdr").Value = "min"
    $ws.Range("N$hdr").Value = "max"
    $ws.Range("O$hdr").Value = "diff"
    $ws.Range("M$hdr").Font().Color = 0
    $ws.Range("N$hdr").Font().Color = 0
    $ws.Range("O$hdr").Font().Color = 0

    # First data row: MIN / MAX / ABS(diff) formulas, styled like the data cells (0.000).
    $ws.Range("M$first").Formula = "=MIN(B$($first):K$($last))"
    $ws.Range("N$first").Formula = "=MAX(B$($first):K$($last))"
    $ws.Range("O$first").Formula = "=ABS(M$($first)-N$($first))"
    $ws.Range("M$first").NumberFormat = "0.000"
    $ws.Range("N$first").NumberFormat = "0.000"

    if ($b.fillDown) {
        for ($r = $first + 1; $r -le $last; $r++) {
            $ws.Range("M$r").NumberFormat = "0.000"
            $ws.Range("N$r").NumberFormat = "0.000"
        }
    }
}

$null = $ws.Range("G4").Select()
